# Generate Report for Archive
#
# Refresh the localization-status report: the file
# "8cae301b-ccda-4441-a461-3e875fe07d2d.md" (row 5 of every sheet) is no
# longer "Ready for handoff" -- it has gone back into translation, so its
# Status columns need to reflect "In Translation" instead, on the
# Overview sheet as well as on each per-locale sheet.

$wb = $excel.ActiveWorkbook

# Overview sheet: zh-cn (E) and de-de (F) status columns for that row.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E5").Value = "In Translation"
$overview.Range("F5").Value = "In Translation"

# zh-cn sheet: Status column (C) for the same file's row.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C5").Value = "In Translation"

# de-de sheet: Status column (C) for the same file's row.
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C5").Value = "In Translation"
